# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets list the same set of events, so both receive the same
# updated values (matching rows differ slightly because "全部类型" has an
# extra row inserted earlier in the sheet).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1329
$ws1.Range("F3").Value  = 1220
$ws1.Range("F4").Value  = 14612
$ws1.Range("F5").Value  = 17704
$ws1.Range("F6").Value  = 146
$ws1.Range("F9").Value  = 217
$ws1.Range("F16").Value = 44
$ws1.Range("F17").Value = 151
$ws1.Range("F19").Value = 1334
$ws1.Range("F22").Value = 64
$ws1.Range("F23").Value = 215
$ws1.Range("F24").Value = 7296
$ws1.Range("F26").Value = 3
$ws1.Range("F27").Value = 40
$ws1.Range("F28").Value = 1176
$ws1.Range("F31").Value = 67
$ws1.Range("F32").Value = 46
$ws1.Range("F33").Value = 142
$ws1.Range("F34").Value = 149
$ws1.Range("F36").Value = 5101

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1329
$ws4.Range("F3").Value  = 1220
$ws4.Range("F4").Value  = 14612
$ws4.Range("F5").Value  = 17704
$ws4.Range("F6").Value  = 146
$ws4.Range("F9").Value  = 217
$ws4.Range("F16").Value = 44
$ws4.Range("F17").Value = 151
$ws4.Range("F19").Value = 1334
$ws4.Range("F23").Value = 64
$ws4.Range("F24").Value = 215
$ws4.Range("F25").Value = 7296
$ws4.Range("F27").Value = 3
$ws4.Range("F28").Value = 40
$ws4.Range("F29").Value = 1176
$ws4.Range("F33").Value = 67
$ws4.Range("F34").Value = 46
$ws4.Range("F35").Value = 142
$ws4.Range("F36").Value = 149
$ws4.Range("F38").Value = 5101
